$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set ExisUnits (column F) to 0 and MaxInvest (column I) to 200 for rows 8-18
for ($r = 8; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 9).Value = 200
}

# Update the selected range to match the author's final selection
$ws.Range("I9:I18").Select()
